$wb = $excel.ActiveWorkbook

# --- Update the absolute path recorded in workbook.xml (x15ac:absPath) ---
$wb.Path = "D:\Jai Mata Dii\DBS_Automation\ExecutionTestData\3\"

# --- DeviceList sheet: remove the first device column (B) and the last two
#     device columns (originally J and K), leaving B..H (7 devices) ---
$ws = $wb.Worksheets.Item("DeviceList")
$ws.Activate()

# Delete column B entirely (shifts C:K left to B:J)
$ws.Columns("B:B").Delete()

# After the shift, the former columns J and K are now I and J - remove them
$ws.Columns("I:J").Delete()

# Individual_ID row is a plain running count (1..7), independent of which
# physical device columns remain - restore the simple sequence
$ws.Cells.Item(9, 2).Value = 1
$ws.Cells.Item(9, 3).Value = 2
$ws.Cells.Item(9, 4).Value = 3
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = 5
$ws.Cells.Item(9, 7).Value = 6
$ws.Cells.Item(9, 8).Value = 7

# Move the selection like the saved workbook shows
$ws.Range("D18").Select()
